$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.944.69"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.793.90"
$ws.Range("E3").Value = "  -0.82%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "602.02"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "163.37"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("D7").Value = "3.791.30"
$ws.Range("E7").Value = "  -0.91%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("E10").Value = "  -2.03%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.88"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +9.33%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.447"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("D15").Value = "4.430.26"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "3.789.93"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "67.942.96"
$ws.Range("E17").Value = "  +0.16%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "18.18"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("E19").Value = "  +1.98%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.02"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.05%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "458.91"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.99%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.47"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.85%  "
$ws.Range("E23").Value = "  -1.51%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.24"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.26%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0000144"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -4.13%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "11.89"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("E28").Value = "  -0.18%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.92"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").Value = "3.939.25"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.59"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -6.93%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("E33").Value = "  -1.77%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "29.03"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("E35").Value = "  +0.09%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.74%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0991"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.16%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.146"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.62%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.81"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("E41").Value = "  -1.98%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +0.12%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "43.73"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.27%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "47.15"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.98%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "152.15"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  -2.15%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.29"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("E50").Value = "  -0.74%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "26.44"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -7.31%  "
